$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("products")

# Replace the full Windows file paths in column D with just the bare file
# names (the "data\images\" prefix is no longer kept in the workbook).
$ws.Cells.Item(2, 4).Value = "controller.png"
$ws.Cells.Item(3, 4).Value = "nail file.jpg"
$ws.Cells.Item(4, 4).Value = "backpack.jpg"
$ws.Cells.Item(5, 4).Value = "pen.jpg"
$ws.Cells.Item(6, 4).Value = "phone.jpg"
$ws.Cells.Item(7, 4).Value = "sunglasses.jpg"

# Move/save the active selection on the sheet to D7 (matches last edited cell).
$ws.Range("D7").Select()
